$d = $word.ActiveDocument

$startPara = 37
$endPara = 55

$targetStart = $d.Paragraphs.Item($startPara).Range.Start
$targetEnd = $d.Paragraphs.Item($endPara).Range.End
$target = $d.Range($targetStart, $targetEnd)

$newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body><w:p>
  <w:r>
    <w:t>Objective and lives</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>The objective of the game is to prevent a meteor from hitting the ground. If that happens then the game is over. Each building has a different amount of meteor hits it can take. The tall and medium buildings can take 3 hits before the 4</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>th</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> meteor hits the ground. The small buildings can only take 2 hits before the 3</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>rd</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> meteor hits the ground. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">If the user-controlled gun is destroyed then the game is over, meaning if even one meteor hits your, building the game is over. </w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>Meteors</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t xml:space="preserve">Meteor 1 has 2 health and goes at an average speed compared to the rest. Meteor 2 (orange) has 3 </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>health</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve"> and travels </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">very </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">slow when compared to the rest. Meteor 3 (yellow) has 1 health and travels fast when compared to the rest. Meteor </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">4 has 8 </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>health</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve"> and travels </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">slowly when compared to the rest. </w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>Guns</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>As the player</w:t>
  </w:r>
  <w:r>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> you only have full control of one gun. You can control that gun using the RIGHT and LEFT arrows and press the SPACE bar to fire. There is a rate of fire control, meaning you can’t blanket the entire screen with your bullets.  Your bullets become faster as you fire at meteors closer to the ground. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">If your gun is destroyed then the game is over. </w:t>
  </w:r>
  <w:r>
    <w:t>The other guns (</w:t>
  </w:r>
  <w:r>
    <w:t>once</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> you have upgraded them) will fire by themselves however they like. </w:t>
  </w:r>
  <w:r>
    <w:t>They</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> aren’t the most accurate. You can also upgrade your own gun.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>Currency</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> and scoring </w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t xml:space="preserve">The </w:t>
  </w:r>
  <w:r>
    <w:t>user-controlled</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> gun gains more fusion</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> cores</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> and points for shooting</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> down a meteor than a computer-controlled</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> gun. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">30 fusion cores for initial upgrade. 60 fusion cores for further upgrade (including </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t xml:space="preserve">user controlled gun). At the end of the game your left over fusion cores are </w:t>
  </w:r>
  <w:r>
    <w:t>multiplied by 10</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> and </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">added to your score. </w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>Upgrades</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t xml:space="preserve">Once you receive enough fusion cores a randomly chosen turret base’s upgrade button will be lit up. Press U to activate the upgrade. You have no control on which turret to perform the upgrade on.  The initially upgraded </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>guns</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t>, including the one the user starts with, does 1 damage per bullet. However the second upgrade does 2 damage.</w:t>
  </w:r>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($newXml)
